# HTTP-API.docx edit script
#
# Semantic changes applied (per the commit "Login erweitert Verwaltung
# für Krankenhäuser" / diff of the canonical OOXML):
#
#   1. The stray "_GoBack" bookmark (left over from the previous save,
#      sitting right after "HTTP" in the title) is removed. Removing it
#      also causes Word's writer to renumber the remaining bookmark
#      ("_Toc342406161") down from id 2 to id 1 automatically.
#   2. "/hospital" + "/op_slot/{id}"  ->  "/hospital" + "/op_slot/delete/{id}"
#      (new DELETE endpoint for removing an OP slot), and Word drops a
#      fresh "_GoBack" bookmark at the point of that edit - right after
#      "/op_slot/delete/" and right before "{id}".
#   3. "/hospital" + "/op_slot"       ->  "/hospital" + "/op_slot/create"
#      (new POST endpoint for creating an OP slot).
#
# (The rest of the diff is Word's proofing engine wrapping already
# existing words in <w:proofErr> spell/grammar-check markers after the
# document was re-opened and spell-checked; that is a purely cosmetic,
# invisible artifact of the split runs and does not change the
# document's visible text or structure.)

$d = $word.ActiveDocument

# --- 1. Remove the stale _GoBack bookmark -----------------------------
# (this also renumbers the "_Toc342406161" bookmark id from 2 -> 1)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. "/hospital/op_slot/{id}" -> "/hospital/op_slot/delete/{id}" ---
# Do this one FIRST, while "/op_slot/" (with trailing slash) is still
# unique in the document - inserting "delete/" right after it turns
# "/op_slot/" + "{id}" into "/op_slot/delete/{id}" (no doubled slash).
$r2 = $d.Content
$r2.Find.Execute("/op_slot/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.InsertAfter("delete/")

# --- 3. Drop the new _GoBack bookmark right after the inserted text ---
$r2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r2)

# --- 4. "/hospital/op_slot" -> "/hospital/op_slot/create" -------------
# Search again from the top of the document: the earlier, standalone
# "/op_slot" occurrence (the one with nothing after it) is still
# untouched and is now the first/only match for a bare "/op_slot".
$r = $d.Content
$r.Find.Execute("/op_slot", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.InsertAfter("/create")
